$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "27.705.98"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.36%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.889.42"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +1.18%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -1.33%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "313.26"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.21%  "

$ws.Range("E6").Value = "  -1.17%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4851"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.26%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3790"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.66%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07331"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.38%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.9190"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.72%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "20.48"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.64%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.07683"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.71%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.895.26"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.45%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.460"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.17%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "6.598"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.37%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "90.89"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.44%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.28%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000008796"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.63%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.15%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "27.744.38"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.51%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "14.52"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.35%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.118"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.22%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "2.151.77"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.75%  "

$ws.Range("E24").Value = "  +0.54%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.912"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.35%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "153.23"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -2.07%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "18.37"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.85%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.114"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +3.92%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "115.83"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.15%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.897"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -1.15%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.08928"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.37%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.151"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -5.41%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.221"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.83%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.7618"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.00%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "4.628"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.50%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.02034"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.16%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.546"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -5.88%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.089"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -4.01%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.05249"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -2.64%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.5454"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -4.10%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.972"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.38%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "6.939"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -1.55%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.1518"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.66%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "8.318"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -2.78%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "109.80"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +4.25%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "10.61"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -1.46%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.4779"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -2.57%  "

$ws.Range("E48").Value = "  -1.22%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.633"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.00%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "67.38"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.28%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.06054"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.84%  "
